$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose refreshed price text looks like a plain number
# (e.g. "1.000", "19.40", "0.05770") must be protected with a Text number
# format before assignment, otherwise Excel auto-converts them to numeric
# values and silently drops significant trailing/leading zeros.
$ws.Range("D5:D11").NumberFormat = "@"
$ws.Range("D13:D20").NumberFormat = "@"
$ws.Range("D22:D23").NumberFormat = "@"
$ws.Range("D25:D51").NumberFormat = "@"

# Column B updates
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("B10").Value = 'Polygon'
$ws.Range("B11").Value = 'Solana'
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("B15").Value = 'TRON'
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("B20").Value = 'Dai'
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("B26").Value = 'Monero'
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("B31").Value = 'Stellar'
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("B36").Value = 'Hedera'
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("B38").Value = 'VeChain'
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("B41").Value = 'PEPE'
$ws.Range("B43").Value = 'Aptos'
$ws.Range("B44").Value = 'MXToken'
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("B47").Value = 'Cronos'
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("B50").Value = 'Quant'
$ws.Range("B51").Value = 'WOONetwork'

# Column C updates
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'

# Column D updates
$ws.Range("D2").Value = '29.280.07'
$ws.Range("D3").Value = '1.931.89'
$ws.Range("D5").Value = '325.54'
$ws.Range("D6").Value = '1.000'
$ws.Range("D7").Value = '0.4623'
$ws.Range("D8").Value = '0.3873'
$ws.Range("D9").Value = '0.07824'
$ws.Range("D10").Value = '0.9715'
$ws.Range("D11").Value = '22.59'
$ws.Range("D12").Value = '1.940.28'
$ws.Range("D13").Value = '7.072'
$ws.Range("D14").Value = '5.771'
$ws.Range("D15").Value = '0.07065'
$ws.Range("D16").Value = '86.74'
$ws.Range("D17").Value = '1.003'
$ws.Range("D18").Value = '0.000009768'
$ws.Range("D19").Value = '17.04'
$ws.Range("D20").Value = '1.002'
$ws.Range("D21").Value = '29.336.06'
$ws.Range("D22").Value = '5.469'
$ws.Range("D23").Value = '11.05'
$ws.Range("D24").Value = '2.154.63'
$ws.Range("D25").Value = '2.092'
$ws.Range("D26").Value = '156.88'
$ws.Range("D27").Value = '19.40'
$ws.Range("D28").Value = '5.758'
$ws.Range("D29").Value = '118.36'
$ws.Range("D30").Value = '1.857'
$ws.Range("D31").Value = '0.09336'
$ws.Range("D32").Value = '0.8634'
$ws.Range("D33").Value = '5.171'
$ws.Range("D34").Value = '1.305'
$ws.Range("D35").Value = '3.077'
$ws.Range("D36").Value = '0.05770'
$ws.Range("D37").Value = '1.154'
$ws.Range("D38").Value = '0.02082'
$ws.Range("D39").Value = '7.665'
$ws.Range("D40").Value = '0.5657'
$ws.Range("D41").Value = '0.000003055'
$ws.Range("D42").Value = '0.1777'
$ws.Range("D43").Value = '9.376'
$ws.Range("D44").Value = '2.718'
$ws.Range("D45").Value = '0.5273'
$ws.Range("D46").Value = '11.44'
$ws.Range("D47").Value = '0.06868'
$ws.Range("D48").Value = '2.078'
$ws.Range("D49").Value = '1.809'
$ws.Range("D50").Value = '111.51'
$ws.Range("D51").Value = '0.2996'

# Column E updates
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  +0.08%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("E22").Value = '  +2.83%  '
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("E28").Value = '  -2.49%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("E35").Value = '  -2.75%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  -1.81%  '
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").Value = '  +53.39%  '
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("E44").Value = '  +6.79%  '
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("E49").Value = '  -1.87%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("E51").Value = '  +1.43%  '

# Restore the original General format on the cells we text-protected above
$ws.Range("D5:D11").NumberFormat = "General"
$ws.Range("D13:D20").NumberFormat = "General"
$ws.Range("D22:D23").NumberFormat = "General"
$ws.Range("D25:D51").NumberFormat = "General"

Write-Host "Applied cryptos update"